# Applies the "npm run build" data-refresh edit to the sample sheet:
#  - L:R columns for rows 2-18 (relatedDoc1-5 / editor_img / editor_bio),
#    previously empty inline strings, now hold the literal text "None".
#  - The editor_social JSON blob in column T (rows 2-18) changes the
#    Python-style `'url': None` entries to quoted `'url': 'None'`.
#  - G9 and G12 ("date") switch from plain text dates to real date
#    serials formatted as YYYY-MM-DD (matching the existing date style
#    already used elsewhere in the sheet, e.g. G2/G5/G10).
#  - A11 ("Status") changes from "Published" to "Draft".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$relatedCols = @("L", "M", "N", "O", "P", "Q", "R")
$socialJson = "[{'icon': 'fab fa-facebook-f', 'url': 'None'}, {'icon': 'fa-brands fa-x-twitter', 'url': 'None'}, {'icon': 'fas fa-link', 'url': 'None'}]"

for ($row = 2; $row -le 18; $row++) {
    foreach ($col in $relatedCols) {
        $ws.Range("$col$row").Value = "None"
    }
    $ws.Range("T$row").Value = $socialJson
}

# Row 9 date: 2024-02-24
$ws.Range("G9").Value = 45346
$ws.Range("G9").NumberFormat = "YYYY-MM-DD"

# Row 12 date: 2024-04-10
$ws.Range("G12").Value = 45392
$ws.Range("G12").NumberFormat = "YYYY-MM-DD"

# Row 11 status: Published -> Draft
$ws.Range("A11").Value = "Draft"
